$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values in rows 1-3
$ws.Range("A1").Value = -0.010969815830924745
$ws.Range("B1").Value = 0.01096981374143109
$ws.Range("A2").Value = 0.0072252819790232072
$ws.Range("B2").Value = -0.0072252841186432014
$ws.Range("A3").Value = 0.027571759699372839
$ws.Range("B3").Value = -0.02757176183540266

# Add new rows 4 and 5
$ws.Range("A4").Value = -0.072933828027747255
$ws.Range("B4").Value = 0.072933826042211131
$ws.Range("A5").Value = 0.030964096205227373
$ws.Range("B5").Value = -0.030964098253871239

# Swap column widths: col A gets old col B width, col B gets old col A width
$ws.Columns.Item(1).ColumnWidth = 13.8
$ws.Columns.Item(2).ColumnWidth = 14.7
